$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B59 was incorrectly stored as text "2"; fix it to be a real numeric value (2)
$ws.Range("B59").Value = 2

# Add new row 60 with the new annotation entry
$ws.Range("A60").Value = "Ruilin"

# B60 must stay a text value "3" (like the rest of that row's raw data),
# not get auto-coerced into a number by COM's type inference.
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = "3"
$ws.Range("B60").ClearFormats()

$ws.Range("C60").Value = "无"
$ws.Range("D60").Value = "FBK"
$ws.Range("E60").Value = "WRI"
$ws.Range("F60").Value = "f6da2ad4-28ad-4a7e-bf94-2041c47bfd2f"
$ws.Range("G60").Value = "rk07ZXZRb_annotated.xlsx"
$ws.Range("H60").Value = "We will add this to the discussion to the paper."
